$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8183
$ws.Range("C3:C8").Value = 7728
$ws.Range("C9:C25").Value = 7660
$ws.Range("C26:C59").Value = 7318
$ws.Range("C213:C252").Value = 7310
